$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '45.294.69'
$ws.Range("E2").Value = '  -0.54%  '

# Row 3
$ws.Range("D3").Value = '2.366.45'
$ws.Range("E3").Value = '  -0.77%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.86'
$ws.Range("E5").Value = '  +0.38%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.42'
$ws.Range("E6").Value = '  -5.21%  '

# Row 7
$ws.Range("E7").Value = '  +0.53%  '

# Row 8
$ws.Range("E8").Value = '  +0.04%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.614'
$ws.Range("E9").Value = '  -2.22%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.96'
$ws.Range("E10").Value = '  -4.68%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0922'
$ws.Range("E11").Value = '  -1.58%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.49'
$ws.Range("E12").Value = '  -2.66%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.110'
$ws.Range("E13").Value = '  +0.26%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.982'
$ws.Range("E14").Value = '  -2.89%  '

# Row 15
$ws.Range("D15").Value = '2.725.91'
$ws.Range("E15").Value = '  -0.77%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.42'
$ws.Range("E16").Value = '  -2.80%  '

# Row 17
$ws.Range("D17").Value = '2.342.79'
$ws.Range("E17").Value = '  -1.28%  '

# Row 18
$ws.Range("D18").Value = '45.190.70'
$ws.Range("E18").Value = '  -0.61%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '15.38'
$ws.Range("E19").Value = '  +13.91%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.29'
$ws.Range("E20").Value = '  -3.72%  '

# Row 21
$ws.Range("E21").Value = '  -1.48%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.61'
$ws.Range("E22").Value = '  +2.58%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.26'
$ws.Range("E23").Value = '  -1.92%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '263.85'
$ws.Range("E24").Value = '  -1.28%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").Value = '  -0.86%  '

# Row 26
$ws.Range("E26").Value = '  +0.04%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.20'
$ws.Range("E27").Value = '  -0.94%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.45'
$ws.Range("E28").Value = '  -4.00%  '

# Row 29
$ws.Range("E29").Value = '  -1.93%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.45'
$ws.Range("E30").Value = '  -1.93%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0951'
$ws.Range("E31").Value = '  -1.78%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.19'
$ws.Range("E32").Value = '  -4.17%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '168.88'
$ws.Range("E33").Value = '  -1.68%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.86'
$ws.Range("E34").Value = '  -4.22%  '

# Row 35
$ws.Range("E35").Value = '  +0.09%  '

# Row 36
$ws.Range("E36").Value = '  -3.70%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.72'
$ws.Range("E37").Value = '  -4.84%  '

# Row 38
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.04'
$ws.Range("E38").Value = '  -0.54%  '

# Row 39
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.03'
$ws.Range("E39").Value = '  -2.83%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.03'
$ws.Range("E40").Value = '  -2.83%  '

# Row 41
$ws.Range("E41").Value = '  -3.03%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.42'
$ws.Range("E42").Value = '  -1.85%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.22'
$ws.Range("E43").Value = '  -2.05%  '

# Row 44
$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.98'
$ws.Range("E44").Value = '  -2.09%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.228'
$ws.Range("E45").Value = '  -4.52%  '

# Row 46
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.859.24'
$ws.Range("E46").Value = '  +13.24%  '

# Row 47
$ws.Range("E47").Value = '  -0.09%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.97'
$ws.Range("E48").Value = '  +3.21%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '84.00'
$ws.Range("E49").Value = '  +6.24%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.08'
$ws.Range("E50").Value = '  -3.52%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.21'
$ws.Range("E51").Value = '  -1.65%  '

Write-Output "Updated cryptos list"